$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 22 de Junio de 2020 a las 11:45"

# Update country statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Muertes hoy, Muertes)
# Row 20: Banglades
$ws.Range("B20").Value = 115786
$ws.Range("C20").Value = 3480
$ws.Range("D20").Value = 46755
$ws.Range("E20").Value = 67529
$ws.Range("G20").Value = 38
$ws.Range("H20").Value = 1502
# Row 32: Indonesia
$ws.Range("B32").Value = 46845
$ws.Range("C32").Value = 954
$ws.Range("D32").Value = 18735
$ws.Range("E32").Value = 25610
$ws.Range("G32").Value = 35
$ws.Range("H32").Value = 2500
# Row 39: Polonia
$ws.Range("B39").Value = 32227
$ws.Range("C39").Value = 296
$ws.Range("E39").Value = 13792
$ws.Range("G39").Value = 3
$ws.Range("H39").Value = 1359
# Row 41: Irak
$ws.Range("B41").Value = 31076
$ws.Range("C41").Value = 1605
$ws.Range("D41").Value = 16408
$ws.Range("E41").Value = 14531
$ws.Range("G41").Value = 6
$ws.Range("H41").Value = 137
# Row 42: Filipinas
$ws.Range("B42").Value = 30868
$ws.Range("D42").Value = 13935
$ws.Range("E42").Value = 15833
$ws.Range("H42").Value = 1100
# Row 43: Oman
$ws.Range("B43").Value = 30052
$ws.Range("D43").Value = 7893
$ws.Range("E43").Value = 20990
$ws.Range("H43").Value = 1169
# Row 51: Israel
$ws.Range("B51").Value = 20869
$ws.Range("C51").Value = 91
$ws.Range("D51").Value = 15714
$ws.Range("E51").Value = 4848
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 307
# Row 68: Marruecos
$ws.Range("B68").Value = 10079
$ws.Range("C68").Value = 102
$ws.Range("D68").Value = 8319
$ws.Range("E68").Value = 1546
# Row 71: Sudan
$ws.Range("B71").Value = 8587
$ws.Range("C71").Value = 15
$ws.Range("D71").Value = 8177
$ws.Range("E71").Value = 289
$ws.Range("H71").Value = 121
# Row 72: Malasia
$ws.Range("B72").Value = 8580
$ws.Range("D72").Value = 3325
$ws.Range("E72").Value = 4734
$ws.Range("H72").Value = 521
# Row 77: Senegal
$ws.Range("B77").Value = 5924
$ws.Range("C77").Value = 98
$ws.Range("D77").Value = 856
$ws.Range("E77").Value = 4933
$ws.Range("G77").Value = 5
$ws.Range("H77").Value = 135
# Row 78: Consejo Danes para los Refugiados
$ws.Range("B78").Value = 5888
$ws.Range("D78").Value = 3919
$ws.Range("E78").Value = 1885
$ws.Range("H78").Value = 84
# Row 105: Estonia
$ws.Range("B105").Value = 1995
$ws.Range("C105").Value = 33
$ws.Range("D105").Value = 1159
$ws.Range("E105").Value = 792
$ws.Range("H105").Value = 44
# Row 106: Albania
$ws.Range("B106").Value = 1981
$ws.Range("D106").Value = 1765
$ws.Range("E106").Value = 147
$ws.Range("H106").Value = 69
# Row 107: Sri Lanka
$ws.Range("D107").Value = 1526
$ws.Range("E107").Value = 413
# Row 118: Eslovenia
$ws.Range("B118").Value = 1521
$ws.Range("C118").Value = 1
$ws.Range("E118").Value = 36
# Row 123: Tunez
$ws.Range("B123").Value = 1159
$ws.Range("C123").Value = 2
$ws.Range("E123").Value = 89
# Row 138: Uganda
$ws.Range("B138").Value = 774
$ws.Range("C138").Value = 4
$ws.Range("D138").Value = 631
$ws.Range("E138").Value = 143
# Row 211: Montserrat
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0
# Row 212: Seychelles
$ws.Range("D212").Value = 10
$ws.Range("H212").Value = 1
# Row 214: Islas Virgenes Britanicas
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
# Row 215: Papua Nueva Guinea
$ws.Range("D215").Value = 7
$ws.Range("H215").Value = 1
